$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 100078720
$ws.Range("B2").Value = 56411
$ws.Range("E2").Value = 100049
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = 'Spillkråka'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = 'Dryocopus martius'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '(Linnaeus, 1758)'
$ws.Range("I2").Value = ""
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = 'äldre spår'
$ws.Range("Q2").Value = 892140.59230468
$ws.Range("R2").Value = 7315543.811931515
$ws.Range("A3").Value = 100078705
$ws.Range("B3").Value = 56540
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 103021
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = 'Talltita'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = 'Poecile montanus'
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = '1'
$ws.Range("K3").Value = ""
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = 'spel/sång'
$ws.Range("Q3").Value = 891870.0593047945
$ws.Range("R3").Value = 7315662.380057276
$ws.Range("A4").Value = 100078718
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = 'äldre spår'
$ws.Range("Q4").Value = 892064.5008710007
$ws.Range("R4").Value = 7315408.872080766
$ws.Range("A5").Value = 100078722
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = 'äldre spår'
$ws.Range("Q5").Value = 892131.3477766616
$ws.Range("R5").Value = 7315498.82523889
$ws.Range("S5").Value = 25
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = '2022-04-13'
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = '2022-04-13'
$ws.Range("AW5").NumberFormat = "@"
$ws.Range("AW5").Value = 'Marie Karlsson'
$ws.Range("AX5").NumberFormat = "@"
$ws.Range("AX5").Value = 'Marie Karlsson'
$ws.Range("A6").Value = 100078725
$ws.Range("B6").Value = 57007
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = 'EN'
$ws.Range("E6").Value = 103042
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = 'Grönfink'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = 'Chloris chloris'
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = '2'
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = 'adult'
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = 'par i lämplig häckbiotop'
$ws.Range("Q6").Value = 891892.9782472391
$ws.Range("R6").Value = 7315677.599053372
$ws.Range("A7").Value = 100078721
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = 'färska spår'
$ws.Range("Q7").Value = 892145.1218588152
$ws.Range("R7").Value = 7315507.794911167
$ws.Range("A8").Value = 100078713
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = 'färska spår'
$ws.Range("Q8").Value = 891337.5114170944
$ws.Range("R8").Value = 7315865.340614381
$ws.Range("A9").Value = 110311519
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = 'födosökande'
$ws.Range("Q9").Value = 891362.1852933455
$ws.Range("R9").Value = 7315800.869760725
$ws.Range("S9").Value = 5
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = '2023-06-23'
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '2023-06-23'
$ws.Range("AW9").NumberFormat = "@"
$ws.Range("AW9").Value = 'Maria Stoltz'
$ws.Range("AX9").NumberFormat = "@"
$ws.Range("AX9").Value = 'Maria Stoltz'
